$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells so numeric-looking strings
# (e.g. "0.7410", "244.40") are not auto-converted to numbers, matching
# the original inline-string cell content.
$priceCells = @("D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.294.50"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.921.33"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").Value = "0.7410"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "244.40"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "0.3143"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "27.33"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("D10").Value = "0.06983"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "0.7749"
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07992"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "1.917.89"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "5.312"
$ws.Range("D15").Value = "91.71"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").Value = "30.355.48"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "14.25"
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("D18").Value = "246.12"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("D19").Value = "5.857"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").Value = "0.000007860"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").Value = "2.195.00"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "1.003"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").Value = "6.678"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").Value = "9.431"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").Value = "165.36"
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("D27").Value = "18.99"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "0.1271"
$ws.Range("E28").Value = "  -5.50%  "
$ws.Range("D29").Value = "2.139"
$ws.Range("E29").Value = "  -7.05%  "
$ws.Range("D30").Value = "1.356"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "1.550"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").Value = "4.358"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "4.081"
$ws.Range("E33").Value = "  -1.82%  "
$ws.Range("D34").Value = "0.05196"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").Value = "1.303"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "0.7499"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "2.775"
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").Value = "0.01946"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").Value = "2.792"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "6.397"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "76.09"
$ws.Range("E41").Value = "  -2.75%  "
$ws.Range("D42").Value = "0.4470"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "1.950"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "0.8383"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("D46").Value = "7.675"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("D47").Value = "101.32"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "9.860"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "37.21"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1223"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "944.66"
$ws.Range("E51").Value = "  -6.22%  "

# Restore default cell style on the Price cells (removes the temporary
# text NumberFormat) so only the cell content differs, not formatting.
foreach ($ref in $priceCells) {
    $ws.Range($ref).Style = "Normal"
}
